# PWR_Board_TestReportTemplate2.xlsx update
# - NoPowerState / NormalOperationState: the HVCAP "TEMP" threshold row (row 7)
#   moves from a +/-0.5 V window to a 200-210 degC window (units relabelled).
# - SPMState: same row 7 TEMP threshold relabelled from V to degC (limits
#   themselves stay 0/5).
# - Quantities: the TEMP quantity's SCALE/OFFSET recalculated (-75 / 200).
# - Also record the new selection/active-sheet state left behind by the edit.

$wb = $excel.ActiveWorkbook

# ---- NoPowerState ---------------------------------------------------------
$wsNoPower = $wb.Worksheets.Item("NoPowerState")
$wsNoPower.Range("B7").Value = 200
$wsNoPower.Range("C7").Value = 210
$wsNoPower.Range("E7").Value = "degC"

# ---- NormalOperationState --------------------------------------------------
$wsNormalOp = $wb.Worksheets.Item("NormalOperationState")
$wsNormalOp.Range("B7").Value = 23.2
$wsNormalOp.Range("C7").Value = 26.8
$wsNormalOp.Range("E7").Value = "degC"

# ---- SPMState ---------------------------------------------------------------
$wsSPM = $wb.Worksheets.Item("SPMState")
$wsSPM.Range("E7").Value = "degC"

# ---- Quantities -------------------------------------------------------------
$wsQuantities = $wb.Worksheets.Item("Quantities")
$wsQuantities.Range("C7").Value = -75
$wsQuantities.Range("D7").Value = 200

# ---- Restore per-sheet selections, then leave NoPowerState as the active tab
$wsReport = $wb.Worksheets.Item("Report")
$wsReport.Activate()
$wsReport.Range("C48").Select()

$wsNormalOp.Activate()
$wsNormalOp.Range("G7").Select()

$wsSPM.Activate()
$wsSPM.Range("F7").Select()

$wsQuantities.Activate()
$wsQuantities.Range("H10").Select()

$wsNoPower.Activate()
$wsNoPower.Range("C7").Select()
